# Applies the cryptos.xlsx price/volume refresh described in the commit
# "Updated cryptos list on Tue Mar 26 06:52:12 UTC 2024 with GitHub Actions".
# Only D (Price) and E (Volume(1h)) cells change; both columns hold plain text
# in the workbook (no numeric cell type), so for any new Price value that COM
# would otherwise auto-coerce into a real number (single-dot decimals such as
# "592.63"), we briefly mark the cell as Text, write the value, then restore the
# original (default) cell style so no stray formatting is introduced.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$defaultStyle = $ws.Range("A1").Style

$ws.Range("D2").Value = "70.442.86"
$ws.Range("E2").Value = "  +4.92%  "
$ws.Range("D3").Value = "3.637.63"
$ws.Range("E3").Value = "  +4.93%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "592.63"
$ws.Range("D5").Style = $defaultStyle
$ws.Range("E5").Value = "  +1.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "194.66"
$ws.Range("D6").Style = $defaultStyle
$ws.Range("E6").Value = "  +3.98%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.647"
$ws.Range("D7").Style = $defaultStyle
$ws.Range("E7").Value = "  +2.27%  "
$ws.Range("D8").Value = "3.629.22"
$ws.Range("E9").Value = "  +0.01%  "
$ws.Range("E10").Value = "  +4.89%  "
$ws.Range("E11").Value = "  +3.69%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "58.49"
$ws.Range("D12").Style = $defaultStyle
$ws.Range("E12").Value = "  +3.97%  "
$ws.Range("E13").Value = "  +4.76%  "
$ws.Range("E14").Value = "  +5.92%  "
$ws.Range("D15").Value = "4.216.21"
$ws.Range("E15").Value = "  +4.83%  "
$ws.Range("E16").Value = "  +6.04%  "
$ws.Range("D17").Value = "3.623.69"
$ws.Range("E17").Value = "  +4.62%  "
$ws.Range("D18").Value = "70.396.84"
$ws.Range("E18").Value = "  +5.08%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.79"
$ws.Range("D19").Style = $defaultStyle
$ws.Range("E19").Value = "  +5.07%  "
$ws.Range("E20").Value = "  +2.13%  "
$ws.Range("E21").Value = "  +5.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "489.31"
$ws.Range("D22").Style = $defaultStyle
$ws.Range("E22").Value = "  +0.40%  "
$ws.Range("E23").Value = "  +13.84%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.35"
$ws.Range("D24").Style = $defaultStyle
$ws.Range("E24").Value = "  +2.12%  "
$ws.Range("E25").Value = "  -0.23%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "91.71"
$ws.Range("D26").Style = $defaultStyle
$ws.Range("E26").Value = "  +2.33%  "
$ws.Range("E27").Value = "  +7.55%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.57"
$ws.Range("D28").Style = $defaultStyle
$ws.Range("E28").Value = "  +5.61%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.64"
$ws.Range("D29").Style = $defaultStyle
$ws.Range("E29").Value = "  +5.79%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "33.07"
$ws.Range("D30").Style = $defaultStyle
$ws.Range("E30").Value = "  +5.32%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.94"
$ws.Range("D31").Style = $defaultStyle
$ws.Range("E31").Value = "  +10.33%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.122"
$ws.Range("D32").Style = $defaultStyle
$ws.Range("E32").Value = "  +9.14%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "630.81"
$ws.Range("D33").Style = $defaultStyle
$ws.Range("E33").Value = "  +5.28%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "12.28"
$ws.Range("D34").Style = $defaultStyle
$ws.Range("E34").Value = "  +4.63%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "65.74"
$ws.Range("D35").Style = $defaultStyle
$ws.Range("E35").Value = "  +2.68%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "41.01"
$ws.Range("D36").Style = $defaultStyle
$ws.Range("E36").Value = "  +12.34%  "
$ws.Range("D38").Value = "0.0₃0829"
$ws.Range("E38").Value = "  +9.43%  "
$ws.Range("E39").Value = "  -1.36%  "
$ws.Range("E40").Value = "  +0.10%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.61"
$ws.Range("D41").Style = $defaultStyle
$ws.Range("E41").Value = "  +1.58%  "
$ws.Range("D42").Value = "3.311.19"
$ws.Range("E42").Value = "  +2.02%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.93"
$ws.Range("D43").Style = $defaultStyle
$ws.Range("E43").Value = "  +15.96%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.19"
$ws.Range("D44").Style = $defaultStyle
$ws.Range("E44").Value = "  +10.05%  "
$ws.Range("E45").Value = "  +5.70%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.91"
$ws.Range("D46").Style = $defaultStyle
$ws.Range("E46").Value = "  +5.32%  "
$ws.Range("E47").Value = "  +2.51%  "
$ws.Range("E48").Value = "  +0.76%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.24"
$ws.Range("D49").Style = $defaultStyle
$ws.Range("E49").Value = "  +5.65%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.33"
$ws.Range("D50").Style = $defaultStyle
$ws.Range("E50").Value = "  +1.46%  "
$ws.Range("E51").Value = "  +0.18%  "
